# Refresh the cryptocurrency ranking table (rows 2-51) on the active sheet
# with newly-scraped values, matching the upstream CSV/CSV-to-xlsx export.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update Price,Volume(1h)
$ws.Range("D2").Value = '51.989.68'
$ws.Range("E2").Value = '  +0.50%  '

# Row 3: update Price,Volume(1h)
$ws.Range("D3").Value = '2.790.68'
$ws.Range("E3").Value = '  -0.95%  '

# Row 4: update Volume(1h)
$ws.Range("E4").Value = '  +0.02%  '

# Row 5: update Price,Volume(1h)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '359.36'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.23%  '

# Row 6: update Price,Volume(1h)
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.69'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.17%  '

# Row 7: update Price,Volume(1h)
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.564'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.32%  '

# Row 8: update Volume(1h)
$ws.Range("E8").Value = '  +0.03%  '

# Row 9: update Volume(1h)
$ws.Range("E9").Value = '  -1.02%  '

# Row 10: update Price,Volume(1h)
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.05'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.25%  '

# Row 11: update Price,Volume(1h)
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0856'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.72%  '

# Row 12: update Volume(1h)
$ws.Range("E12").Value = '  +1.30%  '

# Row 13: update Price,Volume(1h)
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.49'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.83%  '

# Row 14: update Volume(1h)
$ws.Range("E14").Value = '  -2.24%  '

# Row 15: update Price,Volume(1h)
$ws.Range("D15").Value = '3.226.67'
$ws.Range("E15").Value = '  -0.87%  '

# Row 16: update Price,Volume(1h)
$ws.Range("D16").Value = '2.799.71'
$ws.Range("E16").Value = '  -0.87%  '

# Row 17: update Price,Volume(1h)
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.944'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +6.91%  '

# Row 18: update Price,Volume(1h)
$ws.Range("D18").Value = '51.913.77'
$ws.Range("E18").Value = '  +0.68%  '

# Row 19: update Price,Volume(1h)
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.41'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.19%  '

# Row 20: update Price,Volume(1h)
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.13'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.61%  '

# Row 21: update Price,Volume(1h)
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.01'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.72%  '

# Row 22: update Price,Volume(1h)
$ws.Range("D22").Value = '0.0₃0982'
$ws.Range("E22").Value = '  -0.81%  '

# Row 23: update Price,Volume(1h)
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '274.03'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.70%  '

# Row 24: update Price,Volume(1h)
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.27'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.90%  '

# Row 25: update Volume(1h)
$ws.Range("E25").Value = '  +0.19%  '

# Row 26: update Price,Volume(1h)
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.67'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.02%  '

# Row 27: update Volume(1h)
$ws.Range("E27").Value = '  -0.12%  '

# Row 28: update Coin,Link,Price,Volume(1h)
$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.19'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.11%  '

# Row 29: update Coin,Link,Price,Volume(1h)
$ws.Range("B29").Value = 'Kaspa'
$ws.Range("C29").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.145'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.72%  '

# Row 30: update Coin,Link,Price,Volume(1h)
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.21'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.45%  '

# Row 31: update Price,Volume(1h)
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '51.69'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.09%  '

# Row 32: update Coin,Link,Price,Volume(1h)
$ws.Range("B32").Value = 'VeChain'
$ws.Range("C32").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0463'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.93%  '

# Row 33: update Coin,Link,Price,Volume(1h)
$ws.Range("B33").Value = 'InjectiveProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '34.47'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.93%  '

# Row 34: update Price,Volume(1h)
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.74'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.03%  '

# Row 35: update Volume(1h)
$ws.Range("E35").Value = '  +3.30%  '

# Row 36: update Price,Volume(1h)
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.28'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.57%  '

# Row 37: update Volume(1h)
$ws.Range("E37").Value = '  +0.09%  '

# Row 38: update Price,Volume(1h)
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.23'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.52%  '

# Row 39: update Coin,Link,Price,Volume(1h)
$ws.Range("B39").Value = 'Celestia'
$ws.Range("C39").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.12'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.06%  '

# Row 40: update Coin,Link,Price,Volume(1h)
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.64%  '

# Row 41: update Volume(1h)
$ws.Range("E41").Value = '  +1.95%  '

# Row 42: update Volume(1h)
$ws.Range("E42").Value = '  -1.21%  '

# Row 43: update Coin,Link,Price,Volume(1h)
$ws.Range("B43").Value = 'WEMIXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.25'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.72%  '

# Row 44: update Coin,Link,Price,Volume(1h)
$ws.Range("B44").Value = 'Monero'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '122.33'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.03%  '

# Row 45: update Price,Volume(1h)
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '22.20'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.68%  '

# Row 46: update Price,Volume(1h)
$ws.Range("D46").Value = '2.081.72'
$ws.Range("E46").Value = '  +0.26%  '

# Row 47: update Volume(1h)
$ws.Range("E47").Value = '  -1.93%  '

# Row 48: update Price,Volume(1h)
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.22'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.80%  '

# Row 49: update Volume(1h)
$ws.Range("E49").Value = '  +1.33%  '

# Row 50: update Price,Volume(1h)
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.929'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.34%  '

# Row 51: update Price,Volume(1h)
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.93'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.78%  '
